$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header (E1) - "الحصص" (number of class periods / sessions)
$ws.Range("E1").Value = "الحصص"

# Automatic lesson progress tracking: every lesson (rows 2-35) takes 2 sessions
for ($r = 2; $r -le 35; $r++) {
    $ws.Cells.Item($r, 5).Value = 2
}

# Fix navigation: scroll the view and update the current selection
$win = $excel.ActiveWindow
$win.ScrollRow = 21
$win.ScrollColumn = 1
$ws.Range("E2:E35").Select()
